$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits (per sharedStrings diff) ---
# Row 2: FirstName "" -> "rohan"
$ws.Range("B2").Value = "rohan"
# Row 2: Email lini@gmail.com -> tintu@gmail.com
$ws.Range("J2").Value = "tintu@gmail.com"
# Row 2: Hobbies "Reading ,Drawing ,Driving" -> "Reading ,Drawing"
$ws.Range("L2").Value = "Reading ,Drawing"

# Row 3: FirstName "Maya" -> "mini"
$ws.Range("B3").Value = "mini"
# Row 3: Street "dfbdf" -> "abcd"
$ws.Range("H3").Value = "abcd"
# Row 3: Hobbies "Reading ,Writing ," -> "Reading ,Writing"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Row height change (18.75 -> 19.5) for header + data rows ---
$ws.Rows("1:3").RowHeight = 19.5
